{"js": "// Update the division-problem table: each \"a\u00f7b=\" cell text is replaced\n// with a new \"a\u00f7b=\" value, in document order. We walk the table's cells\n// in row-major order and apply the Nth replacement to the Nth non-empty\n// cell, matched against the expected \"before\" text as a safety check.\n// This avoids any ambiguity from duplicate values (e.g. \"92\u00f77=\" appears\n// both as a source and as a target of different replacements).\n\nconst replacements = [\n  \"19\u00f72=|50\u00f79=\",\n  \"31\u00f77=|29\u00f75=\",\n  \"57\u00f72=|85\u00f76=\",\n  \"29\u00f74=|39\u00f72=\",\n  \"62\u00f75=|45\u00f72=\",\n  \"57\u00f78=|30\u00f73=\",\n  \"49\u00f77=|92\u00f77=\",\n  \"54\u00f79=|73\u00f78=\",\n  \"20\u00f76=|72\u00f77=\",\n  \"14\u00f78=|35\u00f73=\",\n  \"59\u00f74=|84\u00f75=\",\n  \"57\u00f79=|13\u00f76=\",\n  \"48\u00f77=|22\u00f72=\",\n  \"73\u00f76=|53\u00f77=\",\n  \"82\u00f73=|33\u00f78=\",\n  \"71\u00f76=|54\u00f73=\",\n  \"92\u00f77=|80\u00f76=\",\n  \"65\u00f73=|85\u00f72=\",\n  \"83\u00f77=|68\u00f79=\",\n  \"40\u00f72=|38\u00f75=\",\n  \"53\u00f79=|62\u00f77=\",\n  \"24\u00f77=|80\u00f72=\",\n  \"91\u00f73=|79\u00f77=\",\n  \"64\u00f73=|68\u00f75=\",\n  \"43\u00f74=|17\u00f75=\",\n].map((pair) => {\n  const [oldText, newText] = pair.split(\"|\");\n  return { oldText, newText };\n});\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Load every row's cells so we can read each cell's current text.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect all cells in row-major (document) order.\nconst allCells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    allCells.push(cell);\n  }\n}\n\n// Load each cell's text.\nfor (const cell of allCells) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\n// Walk the cells in order, applying replacements in sequence whenever a\n// cell's text matches the next expected \"before\" value.\nlet nextReplacement = 0;\nfor (const cell of allCells) {\n  if (nextReplacement >= replacements.length) {\n    break;\n  }\n  const expected = replacements[nextReplacement].oldText;\n  if (cell.body.text === expected) {\n    const newText = replacements[nextReplacement].newText;\n    const cellRange = cell.body.getRange(\"Whole\");\n    cellRange.insertText(newText, Word.InsertLocation.replace);\n    nextReplacement++;\n  }\n}\n\nawait context.sync();\n\nif (nextReplacement !== replacements.length) {\n  throw new Error(\n    \"Only applied \" + nextReplacement + \" of \" + replacements.length + \" replacements\"\n  );\n}\n", "ps1": "# Update the division-problem table: each \"a\u00f7b=\" cell text is replaced\n# with a new \"a\u00f7b=\" value. Cells are addressed directly by (row, column)\n# rather than by text search, since one of the new values (\"92\u00f77=\") also\n# happens to be the pre-existing text of another, later cell \u2014 a plain\n# global find/replace could easily touch the wrong occurrence.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1; Col = 1; OldText = \"19\u00f72=\"; NewText = \"50\u00f79=\" }\n    @{ Row = 1; Col = 2; OldText = \"31\u00f77=\"; NewText = \"29\u00f75=\" }\n    @{ Row = 1; Col = 3; OldText = \"57\u00f72=\"; NewText = \"85\u00f76=\" }\n    @{ Row = 1; Col = 4; OldText = \"29\u00f74=\"; NewText = \"39\u00f72=\" }\n    @{ Row = 1; Col = 5; OldText = \"62\u00f75=\"; NewText = \"45\u00f72=\" }\n    @{ Row = 5; Col = 1; OldText = \"57\u00f78=\"; NewText = \"30\u00f73=\" }\n    @{ Row = 5; Col = 2; OldText = \"49\u00f77=\"; NewText = \"92\u00f77=\" }\n    @{ Row = 5; Col = 3; OldText = \"54\u00f79=\"; NewText = \"73\u00f78=\" }\n    @{ Row = 5; Col = 4; OldText = \"20\u00f76=\"; NewText = \"72\u00f77=\" }\n    @{ Row = 5; Col = 5; OldText = \"14\u00f78=\"; NewText = \"35\u00f73=\" }\n    @{ Row = 9; Col = 1; OldText = \"59\u00f74=\"; NewText = \"84\u00f75=\" }\n    @{ Row = 9; Col = 2; OldText = \"57\u00f79=\"; NewText = \"13\u00f76=\" }\n    @{ Row = 9; Col = 3; OldText = \"48\u00f77=\"; NewText = \"22\u00f72=\" }\n    @{ Row = 9; Col = 4; OldText = \"73\u00f76=\"; NewText = \"53\u00f77=\" }\n    @{ Row = 9; Col = 5; OldText = \"82\u00f73=\"; NewText = \"33\u00f78=\" }\n    @{ Row = 13; Col = 1; OldText = \"71\u00f76=\"; NewText = \"54\u00f73=\" }\n    @{ Row = 13; Col = 2; OldText = \"92\u00f77=\"; NewText = \"80\u00f76=\" }\n    @{ Row = 13; Col = 3; OldText = \"65\u00f73=\"; NewText = \"85\u00f72=\" }\n    @{ Row = 13; Col = 4; OldText = \"83\u00f77=\"; NewText = \"68\u00f79=\" }\n    @{ Row = 13; Col = 5; OldText = \"40\u00f72=\"; NewText = \"38\u00f75=\" }\n    @{ Row = 17; Col = 1; OldText = \"53\u00f79=\"; NewText = \"62\u00f77=\" }\n    @{ Row = 17; Col = 2; OldText = \"24\u00f77=\"; NewText = \"80\u00f72=\" }\n    @{ Row = 17; Col = 3; OldText = \"91\u00f73=\"; NewText = \"79\u00f77=\" }\n    @{ Row = 17; Col = 4; OldText = \"64\u00f73=\"; NewText = \"68\u00f75=\" }\n    @{ Row = 17; Col = 5; OldText = \"43\u00f74=\"; NewText = \"17\u00f75=\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null\n\n    if ($cellRange.Text -ne $edit.OldText) {\n        throw \"Cell ($($edit.Row),$($edit.Col)) expected '$($edit.OldText)' but found '$($cellRange.Text)'\"\n    }\n\n    $cellRange.Text = $edit.NewText\n}\n"}
